$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 567.9
$ws.Range("J12").Value = 619.8570999999999
$ws.Range("L12").Value = 619.8570999999999
$ws.Range("N12").Value = -959.8570999999999
$ws.Range("H32").Value = 6399.2856
$ws.Range("I32").Value = 1996
$ws.Range("J32").Value = 7133.1665
$ws.Range("K32").Value = 1996
$ws.Range("L32").Value = 7133.1665
$ws.Range("M32").Value = -1670
$ws.Range("N32").Value = -7785.1665
$ws.Range("H64").Value = 4373.375
$ws.Range("I64").Value = 4497.8335
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 4497.8335
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -4249.8335
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 4373.375
$ws.Range("I67").Value = 4497.8335
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 4497.8335
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -3639.8335
$ws.Range("N67").Value = -5716
$ws.Range("H80").Value = 732.6591
$ws.Range("I80").Value = 535.2917
$ws.Range("J80").Value = 969.5
$ws.Range("K80").Value = 1605.8751
$ws.Range("L80").Value = 2908.5
$ws.Range("M80").Value = -607.8751
$ws.Range("N80").Value = -4904.5
$ws.Range("H83").Value = 732.6591
$ws.Range("I83").Value = 535.2917
$ws.Range("J83").Value = 969.5
$ws.Range("K83").Value = 4817.6253
$ws.Range("L83").Value = 8725.5
$ws.Range("M83").Value = 174.3747000000003
$ws.Range("N83").Value = -18709.5
$ws.Range("H94").Value = 5083.7144
$ws.Range("I94").Value = 5083.7144
$ws.Range("K94").Value = 5083.7144
$ws.Range("M94").Value = -4632.7144
$ws.Range("H112").Value = 3808.3333
$ws.Range("J112").Value = 3808.3333
$ws.Range("L112").Value = 11424.9999
$ws.Range("N112").Value = -13640.9999
$ws.Range("H113").Value = 6234.5
$ws.Range("I113").Value = 7666.3335
$ws.Range("J113").Value = 5375.4
$ws.Range("K113").Value = 7666.3335
$ws.Range("L113").Value = 5375.4
$ws.Range("M113").Value = -4412.3335
$ws.Range("N113").Value = -11883.4
$ws.Range("H116").Value = 6879.067
$ws.Range("I116").Value = 6957.6
$ws.Range("J116").Value = 6839.8
$ws.Range("K116").Value = 6957.6
$ws.Range("L116").Value = 6839.8
$ws.Range("M116").Value = -3515.6
$ws.Range("N116").Value = -13723.8
$ws.Range("H126").Value = 99988
$ws.Range("J126").Value = 99988
$ws.Range("L126").Value = 99988
$ws.Range("N126").Value = -109868
$ws.Range("H137").Value = 1000.3333
$ws.Range("I137").Value = 839.4
$ws.Range("J137").Value = 1201.5
$ws.Range("K137").Value = 2518.2
$ws.Range("L137").Value = 3604.5
$ws.Range("M137").Value = 31.80000000000018
$ws.Range("N137").Value = -8704.5
$ws.Range("H138").Value = 3872.5
$ws.Range("I138").Value = 2666
$ws.Range("J138").Value = 3997.3103
$ws.Range("K138").Value = 7998
$ws.Range("L138").Value = 11991.9309
$ws.Range("M138").Value = -2858
$ws.Range("N138").Value = -22271.9309

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1941.2593
$ws.Range("I32").Value = 996.7451
$ws.Range("K32").Value = 996.7451
$ws.Range("M32").Value = -709.7451
$ws.Range("H122").Value = 6934.067
$ws.Range("I122").Value = 7671.4165
$ws.Range("J122").Value = 3984.6667
$ws.Range("K122").Value = 23014.2495
$ws.Range("L122").Value = 11954.0001
$ws.Range("M122").Value = -20564.2495
$ws.Range("N122").Value = -16854.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1474.6207
$ws.Range("I20").Value = 1330.3684
$ws.Range("K20").Value = 1330.3684
$ws.Range("M20").Value = -1083.3684
$ws.Range("H105").Value = 2391.5
$ws.Range("I105").Value = 2488.7778
$ws.Range("J105").Value = 2099.6667
$ws.Range("K105").Value = 2488.7778
$ws.Range("L105").Value = 2099.6667
$ws.Range("M105").Value = -741.7777999999998
$ws.Range("N105").Value = -5593.6667
$ws.Range("H107").Value = 2924.625
$ws.Range("I107").Value = 2238.4285
$ws.Range("K107").Value = 2238.4285
$ws.Range("M107").Value = -318.4285
$ws.Range("H134").Value = 1828.8889
$ws.Range("I134").Value = 1837.6923
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 5513.0769
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -2978.0769
$ws.Range("N134").Value = -9870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 916
$ws.Range("I22").Value = 276.30768
$ws.Range("J22").Value = 5074
$ws.Range("K22").Value = 276.30768
$ws.Range("L22").Value = 5074
$ws.Range("M22").Value = 73.69232
$ws.Range("N22").Value = -5774
$ws.Range("H31").Value = 6832.84
$ws.Range("I31").Value = 3627.5
$ws.Range("J31").Value = 7634.175
$ws.Range("K31").Value = 3627.5
$ws.Range("L31").Value = 7634.175
$ws.Range("M31").Value = -3332.5
$ws.Range("N31").Value = -8224.174999999999
$ws.Range("H34").Value = 6832.84
$ws.Range("I34").Value = 3627.5
$ws.Range("J34").Value = 7634.175
$ws.Range("K34").Value = 3627.5
$ws.Range("L34").Value = 7634.175
$ws.Range("M34").Value = -3425.5
$ws.Range("N34").Value = -8038.175
$ws.Range("H69").Value = 15060.333
$ws.Range("I69").Value = 15091
$ws.Range("J69").Value = 14999
$ws.Range("K69").Value = 15091
$ws.Range("L69").Value = 14999
$ws.Range("M69").Value = -14342
$ws.Range("N69").Value = -16497
$ws.Range("H72").Value = 15060.333
$ws.Range("I72").Value = 15091
$ws.Range("J72").Value = 14999
$ws.Range("K72").Value = 45273
$ws.Range("L72").Value = 44997
$ws.Range("M72").Value = -41529
$ws.Range("N72").Value = -52485
$ws.Range("H86").Value = 28574832
$ws.Range("I86").Value = 35717820
$ws.Range("J86").Value = 2873
$ws.Range("K86").Value = 35717820
$ws.Range("L86").Value = 2873
$ws.Range("M86").Value = -35716697
$ws.Range("N86").Value = -5119
$ws.Range("H89").Value = 28574832
$ws.Range("I89").Value = 35717820
$ws.Range("J89").Value = 2873
$ws.Range("K89").Value = 178589100
$ws.Range("L89").Value = 14365
$ws.Range("M89").Value = -178583484
$ws.Range("N89").Value = -25597
$ws.Range("H105").Value = 1179.925
$ws.Range("I105").Value = 1101.08
$ws.Range("J105").Value = 1311.3334
$ws.Range("K105").Value = 1101.08
$ws.Range("L105").Value = 1311.3334
$ws.Range("M105").Value = 645.9200000000001
$ws.Range("N105").Value = -4805.3334
$ws.Range("H107").Value = 730.1042
$ws.Range("I107").Value = 236.09525
$ws.Range("J107").Value = 1114.3334
$ws.Range("K107").Value = 236.09525
$ws.Range("L107").Value = 1114.3334
$ws.Range("M107").Value = 1683.90475
$ws.Range("N107").Value = -4954.3334
$ws.Range("H140").Value = 79396.8
$ws.Range("J140").Value = 79396.8
$ws.Range("L140").Value = 79396.8
$ws.Range("N140").Value = -89756.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2946.5
$ws.Range("I68").Value = 2498.75
$ws.Range("J68").Value = 3125.6
$ws.Range("K68").Value = 7496.25
$ws.Range("L68").Value = 9376.799999999999
$ws.Range("M68").Value = -6685.25
$ws.Range("N68").Value = -10998.8
$ws.Range("H71").Value = 2946.5
$ws.Range("I71").Value = 2498.75
$ws.Range("J71").Value = 3125.6
$ws.Range("K71").Value = 22488.75
$ws.Range("L71").Value = 28130.4
$ws.Range("M71").Value = -18432.75
$ws.Range("N71").Value = -36242.39999999999
$ws.Range("H107").Value = 3172.8667
$ws.Range("I107").Value = 2324.5
$ws.Range("J107").Value = 3481.3635
$ws.Range("K107").Value = 6973.5
$ws.Range("L107").Value = 10444.0905
$ws.Range("M107").Value = -5053.5
$ws.Range("N107").Value = -14284.0905
$ws.Range("H134").Value = 4565
$ws.Range("I134").Value = 4565
$ws.Range("K134").Value = 13695
$ws.Range("M134").Value = -8625
$ws.Range("H141").Value = 5523.3335
$ws.Range("I141").Value = 5523.3335
$ws.Range("K141").Value = 16570.0005
$ws.Range("M141").Value = -11390.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 773.25
$ws.Range("I97").Value = 814.8333
$ws.Range("K97").Value = 814.8333
$ws.Range("M97").Value = -318.8333
$ws.Range("H102").Value = 3240.923
$ws.Range("I102").Value = 3419.3333
$ws.Range("J102").Value = 1100
$ws.Range("K102").Value = 3419.3333
$ws.Range("L102").Value = 1100
$ws.Range("M102").Value = -1797.3333
$ws.Range("N102").Value = -4344
$ws.Range("H126").Value = 7227980.5
$ws.Range("I126").Value = 5465.4443
$ws.Range("K126").Value = 16396.3329
$ws.Range("M126").Value = -13926.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3840.1155
$ws.Range("I61").Value = 3935.1667
$ws.Range("K61").Value = 3935.1667
$ws.Range("M61").Value = -3733.1667
$ws.Range("H100").Value = 4025.5454
$ws.Range("I100").Value = 3766.7036
$ws.Range("J100").Value = 5190.3335
$ws.Range("K100").Value = 3766.7036
$ws.Range("L100").Value = 5190.3335
$ws.Range("M100").Value = -3225.7036
$ws.Range("N100").Value = -6272.3335
$ws.Range("H113").Value = 3840.1155
$ws.Range("I113").Value = 3935.1667
$ws.Range("K113").Value = 3935.1667
$ws.Range("M113").Value = -1765.1667
$ws.Range("H132").Value = 4335.6587
$ws.Range("I132").Value = 3099
$ws.Range("K132").Value = 9297
$ws.Range("M132").Value = -6767

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1699658
$ws.Range("I96").Value = 2667597.5
$ws.Range("J96").Value = 5763.75
$ws.Range("K96").Value = 2667597.5
$ws.Range("L96").Value = 5763.75
$ws.Range("M96").Value = -2666224.5
$ws.Range("N96").Value = -8509.75
$ws.Range("H126").Value = 50005148
$ws.Range("J126").Value = 142861860
$ws.Range("L126").Value = 428585580
$ws.Range("N126").Value = -428590520
$ws.Range("H132").Value = 1050.6774
$ws.Range("I132").Value = 1058.9333
$ws.Range("K132").Value = 3176.7999
$ws.Range("M132").Value = -646.7999
